$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.650.91"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "1.883.81"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'249.30"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.4755"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.2944"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").Value = "'0.06541"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'22.01"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").Value = "'0.07744"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "'96.98"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "'0.7385"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "1.882.13"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "'5.247"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").Value = "'275.43"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "30.623.66"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("E18").Value = "  -3.21%  "
$ws.Range("D19").Value = "'0.000007542"
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("D21").Value = "2.127.54"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'5.352"
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'6.242"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("D25").Value = "'9.233"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").Value = "'164.10"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'18.88"
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "'1.912"
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("D29").Value = "'1.345"
$ws.Range("E29").Value = "  -1.98%  "
$ws.Range("D30").Value = "'0.09735"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").Value = "'1.505"
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").Value = "'4.296"
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").Value = "'4.160"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("D34").Value = "'0.04883"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("D35").Value = "'1.127"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "'0.7003"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").Value = "'2.721"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").Value = "'0.01918"
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("E39").Value = "  +2.36%  "
$ws.Range("D40").Value = "'6.307"
$ws.Range("E40").Value = "  -0.48%  "
$ws.Range("D41").Value = "'75.27"
$ws.Range("E41").Value = "  +6.29%  "
$ws.Range("D42").Value = "'2.038"
$ws.Range("E42").Value = "  +4.62%  "
$ws.Range("D43").Value = "'0.4254"
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D44").Value = "'0.8414"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'102.66"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "'9.394"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("D48").Value = "'7.058"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "'35.63"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "'919.08"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").Value = "'0.05774"
$ws.Range("E51").Value = "  +2.36%  "
